# Applies the "Optuna Attempt (go back with original)" edits to the
# forecast_summary_B0DHWLN5XT_WITH_PO workbook.
#
# Sheet "Forecast Comparison": updates MyForecast (D), Inventory Coverage (H),
#   Stockout Risk (I), Reorder Urgency (J) and Seasonality Index (L) values
#   for rows 2-17.
# Sheet "Summary": updates several aggregate Value cells (B9, B10, B11, B12, B14).

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------

$forecastUpdates = @(
    @{ Row = 2;  D = 84;  H = 9.02; L = 1 },
    @{ Row = 3;  D = 123; H = 5.48; L = 1.03 },
    @{ Row = 4;  D = 116; H = 4.75; L = 1.11 },
    @{ Row = 5;  D = 62;  H = 6.98; L = 1.18 },
    @{ Row = 6;  D = 62;  H = 5.92; L = 0.92 },
    @{ Row = 7;  D = 58;  H = 5.27; L = 0.96 },
    @{ Row = 8;  D = 62;  H = 3.98; I = "Low"; J = "Normal"; L = 1.17 },
    @{ Row = 9;  H = 1.61; I = "Low"; J = "Normal"; L = 0.95 },
    @{ Row = 10; H = 0.63; I = "Low"; L = 0.8100000000000001 },
    @{ Row = 11; L = 1.11 },
    @{ Row = 12; D = 60; L = 1.08 },
    @{ Row = 13; D = 60; L = 0.97 },
    @{ Row = 14; D = 59; L = 1.06 },
    @{ Row = 15; L = 0.84 },
    @{ Row = 16; L = 0.99 },
    @{ Row = 17; L = 0.84 }
)

foreach ($update in $forecastUpdates) {
    $row = $update.Row
    if ($update.ContainsKey("D")) { $wsForecast.Cells.Item($row, 4).Value = $update.D }
    if ($update.ContainsKey("H")) { $wsForecast.Cells.Item($row, 8).Value = $update.H }
    if ($update.ContainsKey("I")) { $wsForecast.Cells.Item($row, 9).Value = $update.I }
    if ($update.ContainsKey("J")) { $wsForecast.Cells.Item($row, 10).Value = $update.J }
    if ($update.ContainsKey("L")) { $wsForecast.Cells.Item($row, 12).Value = $update.L }
}

# --- Summary sheet --------------------------------------------------------------

$wsSummary.Range("B9").Value  = "1405"
$wsSummary.Range("B10").Value = "687"
$wsSummary.Range("B11").Value = "385"
$wsSummary.Range("B12").Value = "123"
$wsSummary.Range("B14").Value = "59"
